# Remove the "networkCrs" parameter row from the scenario_info sheet.
# (commit: "Remove CRS definition from pipeline" - CRS is now read directly
# from the input file instead of being configured here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")

# Row 6 held: networkCrs | EPSG:25832 | The coordinate reference system of the network
# Deleting it shifts every following row up by one (old row 7 -> new row 6, etc.)
[void]$ws.Rows.Item(6).Delete()

# Re-apply the autofilter over the now-smaller range (A1:E25 -> A1:E24)
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:E24").AutoFilter()

# The workbook-level _FilterDatabase defined name also needs to shrink to match
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=scenario_info!`$A`$1:`$E`$24"
    }
}

# scenario_info becomes the active/selected sheet again, with B12 selected
# (previously B11 was selected on the sheet, which is the row right after it
# shifted up to become the sampleSize row's neighbor)
$ws.Activate()
[void]$ws.Range("B12").Select()
